$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = "2019年5月8日23:51:54"
$ws.Range("B42").Value = "周三"
$ws.Range("C42").Value = "dao service"
$ws.Range("D42").Value = "16:00--16:40 & 19:00--21:30"

$ws.Range("C43").Value = "dao service bug修正，未完成"
$ws.Range("D43").Value = "22:40--23:50"

$ws.Range("D43").Select()
